# Refreshes the cryptos price-list snapshot (coinranking.com export)
# to the values captured by the Sun Apr 28 13:51:32 UTC 2024 GitHub
# Actions run: updated Price (D) / Volume(1h) (E) columns, and for rows
# 47-50 the coin rankings reshuffled (Coin/Link/Price/Volume all moved).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '63.820.77'
$ws.Range('E2').Value = '  +0.89%  '

# Row 3
$ws.Range('D3').Value = '3.312.37'
$ws.Range('E3').Value = '  +5.06%  '

# Row 4
$ws.Range('E4').Value = '  +0.00%  '

# Row 5
$ws.Range('E5').Value = '  +1.88%  '

# Row 6
$ws.Range('D6').Value = '''142.39'
$ws.Range('E6').Value = '  +2.76%  '

# Row 7
$ws.Range('E7').Value = '  -0.03%  '

# Row 8
$ws.Range('D8').Value = '3.312.48'
$ws.Range('E8').Value = '  +5.27%  '

# Row 9
$ws.Range('E9').Value = '  +0.66%  '

# Row 10
$ws.Range('E10').Value = '  +2.44%  '

# Row 11
$ws.Range('D11').Value = '''5.48'
$ws.Range('E11').Value = '  +3.67%  '

# Row 12
$ws.Range('D12').Value = '''0.470'
$ws.Range('E12').Value = '  +2.31%  '

# Row 13
$ws.Range('D13').Value = '''0.0000248'
$ws.Range('E13').Value = '  +1.31%  '

# Row 14
$ws.Range('D14').Value = '''34.71'
$ws.Range('E14').Value = '  +1.36%  '

# Row 15
$ws.Range('D15').Value = '3.862.41'
$ws.Range('E15').Value = '  +5.33%  '

# Row 16
$ws.Range('E16').Value = '  +0.04%  '

# Row 17
$ws.Range('D17').Value = '3.313.49'
$ws.Range('E17').Value = '  +5.40%  '

# Row 18
$ws.Range('D18').Value = '63.906.26'
$ws.Range('E18').Value = '  +1.14%  '

# Row 19
$ws.Range('D19').Value = '''6.86'
$ws.Range('E19').Value = '  +2.80%  '

# Row 20
$ws.Range('D20').Value = '''480.35'
$ws.Range('E20').Value = '  +1.18%  '

# Row 21
$ws.Range('D21').Value = '''14.20'
$ws.Range('E21').Value = '  +0.55%  '

# Row 22
$ws.Range('D22').Value = '''0.733'
$ws.Range('E22').Value = '  +4.77%  '

# Row 23
$ws.Range('E23').Value = '  +4.85%  '

# Row 24
$ws.Range('D24').Value = '''13.74'

# Row 25
$ws.Range('D25').Value = '''84.57'
$ws.Range('E25').Value = '  +0.10%  '

# Row 26
$ws.Range('E26').Value = '  +0.11%  '

# Row 27
$ws.Range('E27').Value = '  +1.76%  '

# Row 28
$ws.Range('D28').Value = '''7.33'
$ws.Range('E28').Value = '  +4.58%  '

# Row 29
$ws.Range('E29').Value = '  -0.10%  '

# Row 30
$ws.Range('D30').Value = '''8.13'
$ws.Range('E30').Value = '  +1.05%  '

# Row 31
$ws.Range('E31').Value = '  +2.28%  '

# Row 32
$ws.Range('D32').Value = '''28.79'
$ws.Range('E32').Value = '  +7.09%  '

# Row 33
$ws.Range('E33').Value = '  +0.09%  '

# Row 34
$ws.Range('D34').Value = '''2.55'
$ws.Range('E34').Value = '  +0.77%  '

# Row 35
$ws.Range('E35').Value = '  +3.28%  '

# Row 36
$ws.Range('D36').Value = '''6.00'
$ws.Range('E36').Value = '  +3.34%  '

# Row 37
$ws.Range('D37').Value = '''53.38'
$ws.Range('E37').Value = '  +1.86%  '

# Row 38
$ws.Range('D38').Value = '0.0₃0743'
$ws.Range('E38').Value = '  +6.30%  '

# Row 39
$ws.Range('E39').Value = '  +2.67%  '

# Row 40
$ws.Range('D40').Value = '''433.17'
$ws.Range('E40').Value = '  +2.59%  '

# Row 41
$ws.Range('D41').Value = '3.071.07'
$ws.Range('E41').Value = '  +4.88%  '

# Row 42
$ws.Range('D42').Value = '''2.76'
$ws.Range('E42').Value = '  -0.18%  '

# Row 43
$ws.Range('D43').Value = '''8.35'
$ws.Range('E43').Value = '  +1.29%  '

# Row 44
$ws.Range('E44').Value = '  +1.76%  '

# Row 45
$ws.Range('D45').Value = '''0.265'
$ws.Range('E45').Value = '  +0.85%  '

# Row 46
$ws.Range('D46').Value = '''2.20'
$ws.Range('E46').Value = '  +3.45%  '

# Row 47
$ws.Range('B47').Value = 'Arweave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D47').Value = '''36.40'
$ws.Range('E47').Value = '  +13.74%  '

# Row 48
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').Value = '''26.38'
$ws.Range('E48').Value = '  +3.82%  '

# Row 49
$ws.Range('B49').Value = 'USDe'
$ws.Range('C49').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D49').Value = '''0.999'
$ws.Range('E49').Value = '  +0.02%  '

# Row 50
$ws.Range('B50').Value = 'Monero'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D50').Value = '''127.17'
$ws.Range('E50').Value = '  +5.65%  '

# Row 51
$ws.Range('D51').Value = '''0.114'
$ws.Range('E51').Value = '  +0.99%  '
